# Insert a new "edit" column (B) into the "Main" sheet, shifting the
# existing author/date/yoast_metadesc/excerpt/category/tags/url_path
# columns one place to the right (B:H -> C:I).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# Insert a new column before column B; this shifts B:H -> C:I and
# automatically carries over formatting (including the hyperlink style
# used in column A) to the new column B.
$ws.Columns("B").Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "edit"

# Populate each data row (2-40) with the "edit" hyperlink formula.
for ($r = 2; $r -le 40; $r++) {
    $ws.Range("B$r").Formula = '=HYPERLINK("/wp-admin/post.php?post=&action=edit","edit")'
}
